# Fied import investor access
# - Rename header labels "CF1"/"CF2"/"CF3" (cols S/T/U) to "CF 1"/"CF 2"/"CF 3"
# - Populate new columns T (letters A-H) and U (numbers 100-800) for rows 2-9
# - Update the active selection to U2:U9

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the spacing in the CF1/CF2/CF3 header labels (S1:U1)
$ws.Range("S1").Value = "CF 1"
$ws.Range("T1").Value = "CF 2"
$ws.Range("U1").Value = "CF 3"

# New per-row data in columns T (letter codes) and U (amounts)
$letters = @("A", "B", "C", "D", "E", "F", "G", "H")
$amounts = @(100, 200, 300, 400, 500, 600, 700, 800)

for ($i = 0; $i -lt 8; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 20).Value = $letters[$i]
    $ws.Cells.Item($row, 21).Value = $amounts[$i]
}

$ws.Range("U2:U9").Select()
